$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.229822
$ws.Range("N2").Value = 0.689466
$ws.Range("O2").Value = 0.09226175421862418
$ws.Range("P2").Value = 0.09226175421862419
$ws.Range("Q2").Value = 1.434114100776
$ws.Range("R2").Value = 12.907026906984
$ws.Range("S2").Value = 0.001598186268127886
$ws.Range("T2").Value = 0.001598186268127886

# Row 3
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.4364142651333466
$ws.Range("P3").Value = 0.4364142651333466
$ws.Range("Q3").Value = 6.783611006619999
$ws.Range("R3").Value = 61.05249905957999
$ws.Range("S3").Value = 0.00755970110971989
$ws.Range("T3").Value = 0.007559701109719891

# Row 4
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("M4").Value = 1.174057666666666
$ws.Range("O4").Value = 0.4713239806480292
$ws.Range("P4").Value = 0.4713239806480293
$ws.Range("Q4").Value = 7.326246638227998
$ws.Range("R4").Value = 65.93621974405198
$ws.Range("S4").Value = 0.008164417857545984
$ws.Range("T4").Value = 0.008164417857545988

# Row 5
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.229822
$ws.Range("N5").Value = 0.689466
$ws.Range("O5").Value = 0.09226175421862418
$ws.Range("P5").Value = 0.09226175421862419
$ws.Range("Q5").Value = 79.41880231893799
$ws.Range("R5").Value = 714.7692208704419
$ws.Range("S5").Value = 0.08850484018573561
$ws.Range("T5").Value = 0.08850484018573564

# Row 6
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.4364142651333466
$ws.Range("P6").Value = 0.4364142651333466
$ws.Range("Q6").Value = 375.6648520866016
$ws.Range("S6").Value = 0.4186434034071856
$ws.Range("T6").Value = 0.4186434034071857

# Row 7
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("M7").Value = 1.174057666666666
$ws.Range("O7").Value = 0.4713239806480292
$ws.Range("P7").Value = 0.4713239806480293
$ws.Range("Q7").Value = 405.7150914187222
$ws.Range("R7").Value = 3651.4358227685
$ws.Range("S7").Value = 0.4521315894786877
$ws.Range("T7").Value = 0.4521315894786879

# Row 8
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.229822
$ws.Range("N8").Value = 0.689466
$ws.Range("O8").Value = 0.09226175421862418
$ws.Range("P8").Value = 0.09226175421862419
$ws.Range("Q8").Value = 1.937109577850666
$ws.Range("R8").Value = 17.433986200656
$ws.Range("S8").Value = 0.002158727764760676
$ws.Range("T8").Value = 0.002158727764760677

# Row 9
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.4364142651333466
$ws.Range("P9").Value = 0.4364142651333466
$ws.Range("Q9").Value = 9.162867756635555
$ws.Range("R9").Value = 82.46580980971999
$ws.Range("S9").Value = 0.01021116061644108
$ws.Range("T9").Value = 0.01021116061644109

# Row 10
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("M10").Value = 1.174057666666666
$ws.Range("O10").Value = 0.4713239806480292
$ws.Range("P10").Value = 0.4713239806480293
$ws.Range("Q10").Value = 9.895825251929775
$ws.Range("R10").Value = 89.062427267368
$ws.Range("S10").Value = 0.01102797331179551
$ws.Range("T10").Value = 0.01102797331179551

Write-Host "Update complete"